$d = $word.ActiveDocument

# Update the date line (unique text, use Find/Replace)
$d.Content.Find.Execute("2023-06-30 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-07-01 Saturday", 2) | Out-Null

# Update table cells by position (values are not unique, so address by row/col)
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "99×47=4653"
$t.Cell(1, 2).Range.Text = "93×74=6882"
$t.Cell(1, 3).Range.Text = "82×92=7544"
$t.Cell(1, 4).Range.Text = "67×77=5159"
$t.Cell(1, 5).Range.Text = "55×67=3685"

$t.Cell(2, 1).Range.Text = "50×65=3250"
$t.Cell(2, 2).Range.Text = "24×100=2400"
$t.Cell(2, 3).Range.Text = "38×29=1102"
$t.Cell(2, 4).Range.Text = "72×67=4824"
$t.Cell(2, 5).Range.Text = "85×48=4080"

$t.Cell(3, 1).Range.Text = "68×41=2788"
$t.Cell(3, 2).Range.Text = "24×70=1680"
$t.Cell(3, 3).Range.Text = "64×86=5504"
$t.Cell(3, 4).Range.Text = "40×64=2560"
$t.Cell(3, 5).Range.Text = "41×14=574"

$t.Cell(4, 1).Range.Text = "63×90=5670"
$t.Cell(4, 2).Range.Text = "83×12=996"
$t.Cell(4, 3).Range.Text = "41×55=2255"
$t.Cell(4, 4).Range.Text = "43×16=688"
$t.Cell(4, 5).Range.Text = "89×42=3738"

$t.Cell(5, 1).Range.Text = "64×72=4608"
$t.Cell(5, 2).Range.Text = "16×26=416"
$t.Cell(5, 3).Range.Text = "64×47=3008"
$t.Cell(5, 4).Range.Text = "66×74=4884"
$t.Cell(5, 5).Range.Text = "68×67=4556"

$t.Cell(6, 1).Range.Text = "37×86=3182"
$t.Cell(6, 2).Range.Text = "57×55=3135"
$t.Cell(6, 3).Range.Text = "82×22=1804"
$t.Cell(6, 4).Range.Text = "27×71=1917"
$t.Cell(6, 5).Range.Text = "47×97=4559"

$t.Cell(7, 1).Range.Text = "44×68=2992"
$t.Cell(7, 2).Range.Text = "81×21=1701"
$t.Cell(7, 3).Range.Text = "25×50=1250"
$t.Cell(7, 4).Range.Text = "65×19=1235"
$t.Cell(7, 5).Range.Text = "99×34=3366"

$t.Cell(8, 1).Range.Text = "86×89=7654"
$t.Cell(8, 2).Range.Text = "69×31=2139"
$t.Cell(8, 3).Range.Text = "66×28=1848"
$t.Cell(8, 4).Range.Text = "23×43=989"
$t.Cell(8, 5).Range.Text = "32×16=512"

$t.Cell(9, 1).Range.Text = "29×17=493"
$t.Cell(9, 2).Range.Text = "22×14=308"
$t.Cell(9, 3).Range.Text = "30×71=2130"
$t.Cell(9, 4).Range.Text = "89×67=5963"
$t.Cell(9, 5).Range.Text = "16×24=384"

$t.Cell(10, 1).Range.Text = "72×26=1872"
$t.Cell(10, 2).Range.Text = "56×23=1288"
$t.Cell(10, 3).Range.Text = "18×13=234"
$t.Cell(10, 4).Range.Text = "81×96=7776"
$t.Cell(10, 5).Range.Text = "48×65=3120"

$t.Cell(11, 1).Range.Text = "80×63=5040"
$t.Cell(11, 2).Range.Text = "66×26=1716"
$t.Cell(11, 3).Range.Text = "49×72=3528"
$t.Cell(11, 4).Range.Text = "14×99=1386"
$t.Cell(11, 5).Range.Text = "53×21=1113"

$t.Cell(12, 1).Range.Text = "20×12=240"
$t.Cell(12, 2).Range.Text = "68×90=6120"
$t.Cell(12, 3).Range.Text = "33×60=1980"
$t.Cell(12, 4).Range.Text = "84×85=7140"
$t.Cell(12, 5).Range.Text = "91×20=1820"

$t.Cell(13, 1).Range.Text = "14×66=924"
$t.Cell(13, 2).Range.Text = "55×50=2750"
$t.Cell(13, 3).Range.Text = "35×73=2555"
$t.Cell(13, 4).Range.Text = "89×28=2492"
$t.Cell(13, 5).Range.Text = "59×75=4425"

$t.Cell(14, 1).Range.Text = "86×79=6794"
$t.Cell(14, 2).Range.Text = "95×32=3040"
$t.Cell(14, 3).Range.Text = "13×11=143"
$t.Cell(14, 4).Range.Text = "94×93=8742"
$t.Cell(14, 5).Range.Text = "95×51=4845"

$t.Cell(15, 1).Range.Text = "72×13=936"
$t.Cell(15, 2).Range.Text = "33×28=924"
$t.Cell(15, 3).Range.Text = "43×69=2967"
$t.Cell(15, 4).Range.Text = "82×97=7954"
$t.Cell(15, 5).Range.Text = "32×11=352"

$t.Cell(16, 1).Range.Text = "81×25=2025"
$t.Cell(16, 2).Range.Text = "15×78=1170"
$t.Cell(16, 3).Range.Text = "18×45=810"
$t.Cell(16, 4).Range.Text = "69×53=3657"
$t.Cell(16, 5).Range.Text = "38×36=1368"

$t.Cell(17, 1).Range.Text = "47×33=1551"
$t.Cell(17, 2).Range.Text = "16×56=896"
$t.Cell(17, 3).Range.Text = "19×67=1273"
$t.Cell(17, 4).Range.Text = "41×16=656"
$t.Cell(17, 5).Range.Text = "86×40=3440"

$t.Cell(18, 1).Range.Text = "23×67=1541"
$t.Cell(18, 2).Range.Text = "63×20=1260"
$t.Cell(18, 3).Range.Text = "52×12=624"
$t.Cell(18, 4).Range.Text = "12×97=1164"
$t.Cell(18, 5).Range.Text = "96×47=4512"

$t.Cell(19, 1).Range.Text = "91×73=6643"
$t.Cell(19, 2).Range.Text = "15×81=1215"
$t.Cell(19, 3).Range.Text = "70×71=4970"
$t.Cell(19, 4).Range.Text = "76×70=5320"
$t.Cell(19, 5).Range.Text = "41×71=2911"

$t.Cell(20, 1).Range.Text = "72×35=2520"
$t.Cell(20, 2).Range.Text = "57×85=4845"
$t.Cell(20, 3).Range.Text = "34×13=442"
$t.Cell(20, 4).Range.Text = "12×45=540"
$t.Cell(20, 5).Range.Text = "97×52=5044"
